$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.80449366569519
$ws.Range("B1").Value = 3.759647846221924
$ws.Range("C1").Value = 2.837035417556763
$ws.Range("D1").Value = 0.9159057140350342
$ws.Range("E1").Value = 0.5866419076919556
